# 2020-01-16 공부한 내용 저장
# Adds 4 new rows (13-16) to Sheet2 describing the "find ID" / "find password"
# URL mappings, following the same visual layout/style used by the existing
# rows 11-12 ("로그아웃" section).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 13: 아이디찾기(화면) --------------------------------------------
$ws.Range("A12:H12").Copy($ws.Range("A13:H13"))
$ws.Range("C13").ClearContents()
$ws.Range("H13").ClearContents()

$ws.Range("A13").Value = "아이디찾기(화면)"
$ws.Range("B13").Value = "/member/findID.do"
$ws.Range("D13").Value = "GET"
$ws.Range("E13").Value = "/member/findID.jsp "

# --- Row 14: 아이디찾기(Restful) -----------------------------------------
$ws.Range("A11:H11").Copy($ws.Range("A14:H14"))
$ws.Range("H14").ClearContents()

$ws.Range("A14").Value = "아이디찾기(Restful)"
$ws.Range("B14").Value = "/member/id"
$ws.Range("C14").Value = "tel,birth"
$ws.Range("D14").Value = "GET"
$ws.Range("F14").Value = "MemberSVCImpl"
$ws.Range("G14").Value = "MemberDAOImpl"

# --- Row 15: 비밀번호찾기(화면) ------------------------------------------
$ws.Range("A12:H12").Copy($ws.Range("A15:H15"))
$ws.Range("C15").ClearContents()
$ws.Range("H15").ClearContents()

$ws.Range("A15").Value = "비밀번호찾기(화면)"
$ws.Range("B15").Value = "/member/findPW.do"
$ws.Range("D15").Value = "GET"
$ws.Range("E15").Value = "/member/findPW.jsp "

# --- Row 16: 비밀번호찾기(Restful) ---------------------------------------
$ws.Range("A11:H11").Copy($ws.Range("A16:H16"))
$ws.Range("H16").ClearContents()

$ws.Range("A16").Value = "비밀번호찾기(Restful)"
$ws.Range("B16").Value = "/member/pw"
$ws.Range("C16").Value = "id.tel,birth"
$ws.Range("D16").Value = "GET"
$ws.Range("F16").Value = "MemberSVCImpl"
$ws.Range("G16").Value = "MemberDAOImpl"

# Update selection to match the author's final cursor position.
$ws.Select() | Out-Null
$ws.Range("B8").Select() | Out-Null
